$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of data (2023) appended below the existing table.
$ws.Range("A13").Value = 2023
$ws.Range("B13").Value = 3.5099873380086466
$ws.Range("C13").Value = 8.6650979763673313
$ws.Range("D13").Value = 15.770989815736547

# Match the "0.0" number format already used by the other data rows.
$ws.Range("B13:D13").NumberFormat = "0.0"

# Collapse the old multi-cell selection (A2:D12) back down to the default
# top-left cell now that the table has grown.
$ws.Range("A1").Select()
